# Update hydro progress tracking sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# F3 changed from "In progress" to "Complete"
$ws.Range("F3").Value = "Complete"

# F2/G2 newly filled in with Status/Technician values
$ws.Range("F2").Value = "In progress"
$ws.Range("G2").Value = "Thomas Kosacz"

# Update the active cell selection to F4 (matches saved view state in diff)
$ws.Range("F4").Select()
